$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert " (Please see page 5)" right after "...could be used" and before
#    the following period, inside the reply-to-reviewer-1 paragraph.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("used in the paper should be taken as one possible application where the proposed algorithm could be used")
if (-not $found) {
    throw "Could not locate the target sentence to edit."
}
$matchRange = $find.Parent.Duplicate
$insPoint = $matchRange.Duplicate
$insPoint.Collapse(0)  # wdCollapseEnd
$insPoint.InsertAfter(" (Please see page 5)")

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the document to the empty
#    paragraph that immediately follows the paragraph we just edited.
# ---------------------------------------------------------------------------

# Re-locate the edited paragraph (now containing the inserted text) so we
# can find the paragraph right after it, regardless of exact offsets.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("Indeed our intention to use the context of street light systems")
if (-not $found2) {
    throw "Could not relocate the edited paragraph."
}
$pos = $find2.Parent.Start

$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
        $targetParaIndex = $i
        break
    }
}
if ($targetParaIndex -eq -1) {
    throw "Could not find the paragraph index of the edited paragraph."
}

$followingPara = $d.Paragraphs.Item($targetParaIndex + 1)

# Remove the bookmark from its old location (if present).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-add it, collapsed, inside the following (empty) paragraph. A range that
# spans from just-before to just-after that paragraph's own mark is used so
# the resulting collapsed bookmark resolves to sitting right inside it.
$bmRange = $d.Range($followingPara.Range.Start - 1, $followingPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Edit applied."
